$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "43.582.15"
$ws.Range("E2").Value = "  +1.13%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.381.69"
$ws.Range("E3").Value = "  +3.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "310.26"
$ws.Range("E5").Value = "  -0.02%  "

# Row 6 - Solana
Set-TextValue "D6" "104.95"
$ws.Range("E6").Value = "  +4.00%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -4.35%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.520"
$ws.Range("E9").Value = "  -0.61%  "

# Row 10 - Avalanche
Set-TextValue "D10" "36.24"
$ws.Range("E10").Value = "  +0.97%  "

# Row 11 - OKB
Set-TextValue "D11" "53.44"
$ws.Range("E11").Value = "  +2.53%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0813"
$ws.Range("E12").Value = "  -1.21%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.55%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.02"
$ws.Range("E14").Value = "  -1.36%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.748.09"
$ws.Range("E15").Value = "  +3.54%  "

# Row 16 - Chainlink
$ws.Range("E16").Value = "  +4.63%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.376.69"
$ws.Range("E17").Value = "  +3.34%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +0.79%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "43.517.99"
$ws.Range("E19").Value = "  +1.19%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("E20").Value = "  -4.08%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +3.74%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  -0.38%  "

# Row 23 - Litecoin
Set-TextValue "D23" "68.47"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "241.77"
$ws.Range("E24").Value = "  +0.65%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  +2.58%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +0.35%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.46%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "25.86"
$ws.Range("E28").Value = "  +5.68%  "

# Row 29 - LEO
Set-TextValue "D29" "3.85"
$ws.Range("E29").Value = "  -3.01%  "

# Row 30 - now InjectiveProtocol (was Toncoin)
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "37.04"
$ws.Range("E30").Value = "  -3.63%  "

# Row 31 - now Toncoin (was InjectiveProtocol)
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D31" "2.21"
$ws.Range("E31").Value = "  -4.91%  "

# Row 32 - Cosmos
Set-TextValue "D32" "9.56"
$ws.Range("E32").Value = "  -0.83%  "

# Row 33 - Monero
Set-TextValue "D33" "162.13"
$ws.Range("E33").Value = "  -3.11%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  -0.80%  "

# Row 35 - Celestia
$ws.Range("E35").Value = "  +3.82%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.02%  "

# Row 37 - WEMIXToken
$ws.Range("E37").Value = "  +6.26%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -0.71%  "

# Row 39 - RenderToken
Set-TextValue "D39" "4.70"
$ws.Range("E39").Value = "  +11.00%  "

# Row 40 - Hedera
Set-TextValue "D40" "0.0744"
$ws.Range("E40").Value = "  +0.47%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +6.55%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -1.15%  "

# Row 43 - Stellar
$ws.Range("E43").Value = "  -1.76%  "

# Row 44 - ApeXProtocol
Set-TextValue "D44" "2.66"
$ws.Range("E44").Value = "  +16.04%  "

# Row 45 - now Maker (was EnergySwap)
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D45" "2.032.79"
$ws.Range("E45").Value = "  +2.94%  "

# Row 46 - now EnergySwap (was Maker)
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "19.70"
$ws.Range("E46").Value = "  +3.45%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +0.32%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  +4.04%  "

# Row 49 - FraxShare
Set-TextValue "D49" "10.60"
$ws.Range("E49").Value = "  +7.59%  "

# Row 50 - MultiversX
Set-TextValue "D50" "58.01"
$ws.Range("E50").Value = "  +3.85%  "

# Row 51 - HuobiToken
$ws.Range("E51").Value = "  +0.95%  "
